# Login & Logout set and remove cookie
# Login & Logout set and remove cookie for email
#
# The "Backlog" J-column (column J on Sheet1) gets its text content shifted
# up by one row (J7's old note "Drop cookie using ngcookies - $cookieStore"
# is dropped entirely), and the final entry is replaced with "Encrypt cookie".
# Two other now-unused notes ("Log out - implement" and
# "Remove email adress from URLS - cookies") are removed from the shared
# strings as a consequence of no longer being referenced anywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the J-column cells whose content is being removed/relocated.
$ws.Range("J7").Clear()
$ws.Range("J9").Clear()
$ws.Range("J11").Clear()
$ws.Range("J13").Clear()
$ws.Range("J16").Clear()
$ws.Range("J18").Clear()
$ws.Range("J22").Clear()
$ws.Range("J24").Clear()
$ws.Range("J27").Clear()

# Write the shifted-up content (each row now holds what used to be one row
# below it), keeping existing cell formatting (J8 retains its wrap-text
# style) intact.
$ws.Range("J8").Value = "Limit favourites drop down to show first 5"
$ws.Range("J10").Value = "Add Exception handling - add to db"
$ws.Range("J12").Value = "Trending Now"
$ws.Range("J14").Value = "Check Thumbnail pics for cars - and change ones not adequate"
$ws.Range("J15").Value = "Add missing CarInfo entries"
$ws.Range("J17").Value = "Add Car List as Grid View"
$ws.Range("J20").Value = "1) ORDERS - Add rental orders to the Db"
$ws.Range("J21").Value = "CarInfo - add for each individual car in list"
$ws.Range("J23").Value = "Create 'My account page' - switch on/off Recommended cars"
$ws.Range("J25").Value = "Encrypt cookie"
